# "Metro 52 Noord = Zuid" — flip the timetable to the Noord -> Zuid direction
# (previously the sheet listed the Station Zuid -> Noord run).
# Station-name cells keep their text; the "HH:MM + 1" text times are replaced
# by real time-of-day serial numbers formatted as h:mm (same as the existing
# arrival-time rows already on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Station order reverses: Noord, Noorderpark, Centraal Station, Rokin,
# Vijzelgracht, De Pijp, Europaplein, Station Zuid
$ws.Cells.Item(1, 1).Value  = "Noord"
$ws.Cells.Item(5, 1).Value  = "Noorderpark"
$ws.Cells.Item(9, 1).Value  = "Centraal Station"
$ws.Cells.Item(13, 1).Value = "Rokin"
$ws.Cells.Item(17, 1).Value = "Vijzelgracht"
$ws.Cells.Item(21, 1).Value = "De Pijp"
$ws.Cells.Item(25, 1).Value = "Europaplein"
$ws.Cells.Item(29, 1).Value = "Station Zuid"

# The three departure-time rows under each station become numeric times
# (serial day fractions) instead of literal "HH:MM + 1" text.
$times = [ordered]@{
    2  = 0.7104166666666667
    3  = 0.71388888888888891
    4  = 0.71736111111111101
    6  = 0.71180555555555547
    7  = 0.71527777777777779
    8  = 0.71875
    10 = 0.71319444444444446
    11 = 0.71666666666666667
    12 = 0.72013888888888899
    14 = 0.71458333333333324
    15 = 0.71805555555555556
    16 = 0.72152777777777777
    18 = 0.71597222222222223
    19 = 0.71944444444444444
    20 = 0.72291666666666676
    22 = 0.71736111111111101
    23 = 0.72083333333333333
    24 = 0.72430555555555554
    26 = 0.71875
    27 = 0.72222222222222221
    28 = 0.72569444444444453
    30 = 0.72013888888888899
    31 = 0.72361111111111109
    32 = 0.7270833333333333
}

foreach ($r in $times.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $times[$r]
    $cell.NumberFormat = "h:mm"
}

# View state left by the author: scrolled a bit down, selection on B14.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
[void]$ws.Range("B14").Select()

# Page setup for printing: A4, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "Rewrote timetable to Noord -> Zuid direction (8 stations, 24 times), updated view/page setup."
